$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80").Value = "Login with valid username and password"
$ws.Range("B80").Value = "PASSED"
$ws.Range("C80").Value = "chrome"

$ws.Range("A81").Value = "Create Country"
$ws.Range("B81").Value = "PASSED"
$ws.Range("C81").Value = "chrome"
